# Scheduled-runner update: refresh market-price-derived profit figures
# (currentAveragePrice / LevePrice / LeveProfit columns H-N) across the
# Exodus_Profits sheets, per the latest market data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 179.4
$ws.Range("I33").Value = 178.375
$ws.Range("K33").Value = 178.375
$ws.Range("M33").Value = 50.625

$ws.Range("H64").Value = 4952.2856
$ws.Range("J64").Value = 6999.778
$ws.Range("L64").Value = 6999.778
$ws.Range("N64").Value = -7495.778

$ws.Range("H67").Value = 4952.2856
$ws.Range("J67").Value = 6999.778
$ws.Range("L67").Value = 6999.778
$ws.Range("N67").Value = -8715.778

$ws.Range("H132").Value = 1488.159
$ws.Range("I132").Value = 1487.619
$ws.Range("K132").Value = 4462.857
$ws.Range("M132").Value = -1932.857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5282.085
$ws.Range("I32").Value = 2330.415
$ws.Range("K32").Value = 2330.415
$ws.Range("M32").Value = -2043.415

$ws.Range("H61").Value = 85979.5
$ws.Range("I61").Value = 2552.75
$ws.Range("K61").Value = 2552.75
$ws.Range("M61").Value = -2340.75

$ws.Range("H63").Value = 3423.7273
$ws.Range("I63").Value = 2571.2856
$ws.Range("J63").Value = 4915.5
$ws.Range("K63").Value = 2571.2856
$ws.Range("L63").Value = 4915.5
$ws.Range("M63").Value = -1885.2856
$ws.Range("N63").Value = -6287.5

$ws.Range("H66").Value = 3423.7273
$ws.Range("I66").Value = 2571.2856
$ws.Range("J66").Value = 4915.5
$ws.Range("K66").Value = 12856.428
$ws.Range("L66").Value = 24577.5
$ws.Range("M66").Value = -9424.428
$ws.Range("N66").Value = -31441.5

$ws.Range("H132").Value = 3622.0244
$ws.Range("I132").Value = 3111.0908
$ws.Range("K132").Value = 9333.2724
$ws.Range("M132").Value = -6803.2724

$ws.Range("H136").Value = 85979.5
$ws.Range("I136").Value = 2552.75
$ws.Range("K136").Value = 7658.25
$ws.Range("M136").Value = -5108.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 36999
$ws.Range("J81").Value = 36999
$ws.Range("L81").Value = 36999
$ws.Range("N81").Value = -39121

$ws.Range("H84").Value = 36999
$ws.Range("J84").Value = 36999
$ws.Range("L84").Value = 110997
$ws.Range("N84").Value = -121605

$ws.Range("H105").Value = 145576.72
$ws.Range("I105").Value = 202517.8
$ws.Range("K105").Value = 202517.8
$ws.Range("M105").Value = -200770.8

$ws.Range("H123").Value = 50000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 50000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 50000
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -59800

$ws.Range("H134").Value = 3918.12
$ws.Range("I134").Value = 2403.4102
$ws.Range("K134").Value = 7210.230599999999
$ws.Range("M134").Value = -4675.230599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3767.0952
$ws.Range("I58").Value = 3928
$ws.Range("K58").Value = 3928
$ws.Range("M58").Value = -3725

$ws.Range("H70").Value = 44166.668
$ws.Range("J70").Value = 44166.668
$ws.Range("L70").Value = 44166.668
$ws.Range("N70").Value = -44796.668

$ws.Range("H73").Value = 44166.668
$ws.Range("J73").Value = 44166.668
$ws.Range("L73").Value = 44166.668
$ws.Range("N73").Value = -46350.668

$ws.Range("H122").Value = 5033.385
$ws.Range("I122").Value = 5030
$ws.Range("J122").Value = 5037.3335
$ws.Range("K122").Value = 15090
$ws.Range("L122").Value = 15112.0005
$ws.Range("M122").Value = -12640
$ws.Range("N122").Value = -20012.0005

$ws.Range("H136").Value = 3767.0952
$ws.Range("I136").Value = 3928
$ws.Range("K136").Value = 11784
$ws.Range("M136").Value = -9234

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3417038
$ws.Range("J4").Value = 2857485.8
$ws.Range("L4").Value = 8572457.399999999
$ws.Range("N4").Value = -8572681.399999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 51860.5
$ws.Range("I80").Value = 51860.5
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 51860.5
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -50862.5
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 51860.5
$ws.Range("I83").Value = 51860.5
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 259302.5
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -254310.5
$ws.Range("N83").ClearContents()

$ws.Range("H132").Value = 4062.5
$ws.Range("I132").Value = 4000
$ws.Range("J132").Value = 4250
$ws.Range("K132").Value = 12000
$ws.Range("L132").Value = 12750
$ws.Range("M132").Value = -9470
$ws.Range("N132").Value = -17810

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9916.736999999999
$ws.Range("J7").Value = 4919
$ws.Range("L7").Value = 4919
$ws.Range("N7").Value = -5143

$ws.Range("H74").Value = 44130
$ws.Range("I74").Value = 24999.5
$ws.Range("K74").Value = 24999.5
$ws.Range("M74").Value = -24001.5

$ws.Range("H77").Value = 44130
$ws.Range("I77").Value = 24999.5
$ws.Range("K77").Value = 74998.5
$ws.Range("M77").Value = -70006.5

$ws.Range("H126").Value = 9916.736999999999
$ws.Range("J126").Value = 4919
$ws.Range("L126").Value = 14757
$ws.Range("N126").Value = -19697

$ws.Range("H132").Value = 4104.9165
$ws.Range("I132").Value = 3013.0715
$ws.Range("J132").Value = 5633.5
$ws.Range("K132").Value = 9039.2145
$ws.Range("L132").Value = 16900.5
$ws.Range("M132").Value = -6509.2145
$ws.Range("N132").Value = -21960.5

$ws.Range("H136").Value = 4112.2607
$ws.Range("I136").Value = 4963.5
$ws.Range("K136").Value = 14890.5
$ws.Range("M136").Value = -12340.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 4333
$ws.Range("J3").Value = 4333
$ws.Range("L3").Value = 4333
$ws.Range("N3").Value = -4561

$ws.Range("H96").Value = 5850963.5
$ws.Range("J96").Value = 8774028
$ws.Range("L96").Value = 8774028
$ws.Range("N96").Value = -8776774

$ws.Range("H113").Value = 9853.714
$ws.Range("I113").Value = 7795.4
$ws.Range("K113").Value = 23386.2
$ws.Range("M113").Value = -21216.2

$ws.Range("H122").Value = 3674
$ws.Range("I122").Value = 3830.8572
$ws.Range("J122").Value = 3125
$ws.Range("K122").Value = 11492.5716
$ws.Range("L122").Value = 9375
$ws.Range("M122").Value = -9042.571599999999
$ws.Range("N122").Value = -14275
